$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 11:42"

# Full province data set (Ciudad, Casos totales, Casos activos, Recuperados, Muertes)
# re-sorted by Casos totales (descending) after incorporating updated figures for
# Valencia/Valencia, Alacant/Alicante and Castello/Castellon.
$data = @(
    @("Madrid", 14597, 3031, 9741, 1825),
    @("Cataluña", 11592, 1697, 9223, 672),
    @("Castilla-La Mancha", 2780, 71, 2446, 263),
    @("Valencia/Valencia", 1788, 36, 1680, 72),
    @("Bizkaia/Vizcaya", 1501, 466, 1447, 54),
    @("Araba/Alava", 1207, 466, 1126, 81),
    @("Navarra", 1197, 23, 1141, 33),
    @("Alacant/Alicante", 1039, 17, 941, 81),
    @("La Rioja", 928, 43, 848, 37),
    @("Aragon", 907, 29, 838, 40),
    @("Ciudad Real", 885, 71, 810, 67),
    @("Malaga", 819, 48, 731, 40),
    @("Asturias", 779, 35, 719, 25),
    @("Albacete", 752, 71, 684, 60),
    @("Toledo", 752, 71, 660, 65),
    @("A Coruña", 734, 25, 715, 19),
    @("Pontevedra", 602, 25, 598, 4),
    @("Murcia", 596, 9, 579, 8),
    @("Granada", 579, 1, 553, 25),
    @("Gipuzkoa/Guipuzcoa", 563, 466, 543, 20),
    @("Tenerife", 539, 15, 519, 24),
    @("Sevilla", 535, 6, 511, 18),
    @("Salamanca", 533, 57, 430, 46),
    @("Cantabria", 510, 12, 484, 14),
    @("Valladolid", 501, 36, 444, 21),
    @("Caceres", 485, 3, 447, 35),
    @("Burgos", 485, 55, 402, 28),
    @("Leon", 438, 31, 376, 31),
    @("Guadalajara", 404, 71, 353, 49),
    @("Castello/Castellon", 373, 4, 355, 14),
    @("Zaragoza", 329, 0, 315, 14),
    @("Jaen", 316, 5, 297, 14),
    @("Segovia", 300, 47, 219, 34),
    @("Cordoba", 291, 0, 285, 6),
    @("Cadiz", 278, 4, 270, 4),
    @("Soria", 263, 19, 230, 14),
    @("Badajoz", 257, 5, 248, 4),
    @("Avila", 226, 36, 169, 21),
    @("Mallorca", 210, 18, 194, 12),
    @("Ourense", 189, 25, 186, 3),
    @("Cuenca", 172, 71, 142, 22),
    @("Gran Canaria", 171, 15, 169, 24),
    @("Almeria", 115, 5, 105, 5),
    @("Lugo", 103, 25, 99, 4),
    @("Zamora", 103, 16, 79, 8),
    @("Palencia", 91, 11, 77, 3),
    @("Huelva", 77, 2, 74, 1),
    @("Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena", 58, 0, 58, 3),
    @("Teruel", 47, 0, 45, 2),
    @("Melilla", 38, 0, 37, 1),
    @("Huesca", 37, 0, 37, 0),
    @("La Palma", 33, 15, 33, 24),
    @("Ibiza", 21, 18, 20, 1),
    @("Fuerteventura", 20, 15, 20, 24),
    @("Lanzarote", 17, 15, 17, 24),
    @("Ceuta", 16, 0, 16, 0),
    @("Menorca", 15, 18, 13, 0),
    @("Arroyo de la Luz", 7, 0, 7, 0),
    @("La Gomera", 4, 15, 2, 24),
    @("El Hierro", 3, 15, 3, 24),
    @("Formentera", 0, 10, 0, 8)
)

$row = 4
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}
